# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
#
# The "Date" column (BF) on the sheet stores the game date as a literal
# text string. It was previously written as "5-13-2012-13" (a mangled
# "<month>-<day>-<season>" string); it should instead be the ISO-style
# "2013-05-13" (the actual game date).
#
# NOTE: Simply assigning a string like "2013-05-13" straight to
# Range.Value causes Excel to auto-recognise it as a date and silently
# convert the cell to a date serial number, which is not what we want -
# we need the literal text preserved. Routing the text through a
# TEXT-producing formula on a scratch cell and pasting *values only*
# avoids that automatic date coercion while keeping the target cell's
# number format untouched (General, just like the other text cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDate = "5-13-2012-13"
$newDate = "2013-05-13"

# Scratch cell well outside the used range, used only to generate a
# guaranteed-text value for the paste-special below.
$helper = $ws.Cells.Item(100, 100)

for ($row = 2; $row -le 31; $row++) {
    $target = $ws.Cells.Item($row, 58)  # column BF
    if ($target.Value2 -eq $oldDate) {
        $helper.Formula = "=""$newDate"""
        $helper.Copy()
        $target.PasteSpecial(-4163)  # xlPasteValues
    }
}

$excel.CutCopyMode = $false
$helper.ClearContents()
